$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete rows 5-7 (data consolidated into rows 2-4 with updated TPM values)
$ws.Rows("5:7").Delete()

# Row 2: MuSCs / Cdh1 / Igf1r -> ECs
$ws.Range("A2").Value = "MuSCs"
$ws.Range("B2").Value = "Cdh1"
$ws.Range("C2").Value = "Igf1r"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.9477166666666667
$ws.Range("H2").Value = 2.84315
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 14.129345
$ws.Range("N2").Value = 42.388035
$ws.Range("O2").Value = 0.3414817166893976
$ws.Range("P2").Value = 0.3414817166893976
$ws.Range("Q2").Value = 13.39061574558333
$ws.Range("R2").Value = 120.51554171025
$ws.Range("S2").Value = 0.3414817166893976
$ws.Range("T2").Value = 0.3414817166893976

# Row 3: MuSCs / Cdh1 / Igf1r -> FAPs
$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Cdh1"
$ws.Range("C3").Value = "Igf1r"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.9477166666666667
$ws.Range("H3").Value = 2.84315
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 14.70328633333333
$ws.Range("N3").Value = 44.109859
$ws.Range("O3").Value = 0.3553528814026711
$ws.Range("P3").Value = 0.3553528814026711
$ws.Range("Q3").Value = 13.93454951287222
$ws.Range("R3").Value = 125.41094561585
$ws.Range("S3").Value = 0.3553528814026711
$ws.Range("T3").Value = 0.3553528814026711

# Row 4: MuSCs / Cdh1 / Igf1r -> MuSCs
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Cdh1"
$ws.Range("C4").Value = "Igf1r"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.9477166666666667
$ws.Range("H4").Value = 2.84315
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 12.543947
$ws.Range("N4").Value = 37.631841
$ws.Range("O4").Value = 0.3031654019079313
$ws.Range("P4").Value = 0.3031654019079312
$ws.Range("Q4").Value = 11.88810763768333
$ws.Range("R4").Value = 106.99296873915
$ws.Range("S4").Value = 0.3031654019079313
$ws.Range("T4").Value = 0.3031654019079312
